# Update countries & provincias Spain
# Applies the 22-Oct-2020 15:33 data refresh to the "Pais" sheet:
#  - Paises Bajos overtakes Belgica (rows 29/30 swap label + stats)
#  - Serbia overtakes Bosnia y Herzegovina (rows 81/82 swap label + stats)
#  - Groenlandia overtakes San Pedro y Miquelon (rows 214/215 swap label + stats)
#  - Refreshed case counters for several other countries
#  - Title timestamp bumped from 14:16 to 15:33

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / timestamp
$ws.Range("A1").Value = "Datos actualizados a 22 de Octubre de 2020 a las 15:33"

# Row 19
$ws.Range("B19").Value = 442164
$ws.Range("C19").Value = 3899
$ws.Range("D19").Value = 371826
$ws.Range("E19").Value = 59873
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = 10465

# Row 25
$ws.Range("B25").Value = 343774
$ws.Range("C25").Value = 401
$ws.Range("D25").Value = 330181
$ws.Range("E25").Value = 8343
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 5250

# Row 29 - becomes Paises Bajos
$ws.Range("A29").Value = "Paises Bajos"
$ws.Range("B29").Value = 262405
$ws.Range("C29").Value = 9271
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("G29").Value = 46
$ws.Range("H29").Value = 6919

# Row 30 - becomes Belgica
$ws.Range("A30").Value = "Belgica"
$ws.Range("B30").Value = 253386
$ws.Range("C30").Value = 13227
$ws.Range("D30").Value = 21717
$ws.Range("E30").Value = 221130
$ws.Range("G30").Value = 50
$ws.Range("H30").Value = 10539

# Row 39
$ws.Range("B39").Value = 130462
$ws.Range("C39").Value = 252
$ws.Range("D39").Value = 127328
$ws.Range("E39").Value = 2906
$ws.Range("G39").Value = 3
$ws.Range("H39").Value = 228

# Row 43
$ws.Range("B43").Value = 119420
$ws.Range("C43").Value = 889
$ws.Range("D43").Value = 110714
$ws.Range("E43").Value = 7976
$ws.Range("G43").Value = 9
$ws.Range("H43").Value = 730

# Row 63
$ws.Range("B63").Value = 64439
$ws.Range("C63").Value = 429
$ws.Range("D63").Value = 61658
$ws.Range("E63").Value = 2241
$ws.Range("G63").Value = 6
$ws.Range("H63").Value = 540

# Row 81 - becomes Serbia
$ws.Range("A81").Value = "Serbia"
$ws.Range("B81").Value = 37536
$ws.Range("C81").Value = 416
$ws.Range("D81").Value = 31536
$ws.Range("E81").Value = 5217
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 783

# Row 82 - becomes Bosnia y Herzegovina
$ws.Range("A82").Value = "Bosnia y Herzegovina"
$ws.Range("B82").Value = 37314
$ws.Range("C82").Value = 999
$ws.Range("D82").Value = 25989
$ws.Range("E82").Value = 10274
$ws.Range("G82").Value = 20
$ws.Range("H82").Value = 1051

# Row 124
$ws.Range("E124").Value = 2453
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 14

# Row 163
$ws.Range("B163").Value = 1923
$ws.Range("C163").Value = 5
$ws.Range("E163").Value = 919

# Row 214 - becomes Groenlandia
$ws.Range("A214").Value = "Groenlandia"
$ws.Range("B214").Value = 17
$ws.Range("C214").Value = 1
$ws.Range("D214").Value = 16
$ws.Range("E214").Value = 1

# Row 215 - becomes San Pedro y Miquelon
$ws.Range("A215").Value = "San Pedro y Miquelon"
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 4
